# productos.xlsx - "termino del proyecto con readme"
# Update row 2 (id 1 -> 8503, new media code + image url) and remove the
# now-superseded row 3, matching the final state of the productos sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 ---------------------------------------------------------
# A2: numeric id
$ws.Range("A2").Value = 8503

# C2: plain text url (not numeric-looking, assigns as text directly)
$ws.Range("C2").Value = "https://i.ibb.co/5jcC71J/M105.jpg"

# B2: the new value "150634" looks like a number, so a naive .Value=
# assignment would be stored as a numeric cell. Build it as a text formula
# in a scratch cell (row 3, which gets removed below anyway) and paste only
# the resulting value into B2, so it lands as plain text with no leftover
# number-format/quote-prefix style.
$ws.Range("Z3").Formula = '="150634"'
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4163) | Out-Null  # xlPasteValues
$excel.CutCopyMode = 0

# --- Remove the old row 3 (its data now lives in row 2) ------------------
$ws.Range("A3:XFD1371").EntireRow.Select() | Out-Null
$ws.Range("A3:XFD1371").EntireRow.Delete() | Out-Null
